# daily auto push: 2026-01-31 06:55 UTC
# Insert a new data row just above the existing row 729 (which shifts the
# 2026/12/29 ... 2027/01/05 block down by one row) and populate it with the
# new entry: 2026/01/31 (土), time 13, ranking 22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(729).Insert()

$ws.Cells.Item(729, 1).NumberFormat = "@"
$ws.Cells.Item(729, 1).Value = "2026/01/31"
$ws.Cells.Item(729, 2).Value = "土"
$ws.Cells.Item(729, 3).Value = 13
$ws.Cells.Item(729, 4).Value = 22
